# Updated cryptos list values (Price and Volume(1h)) per source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.752.84"
$ws.Range("E2").Value = "  -0.22%  "
$ws.Range("D3").Value = "2.249.78"
$ws.Range("E3").Value = "  +0.57%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "112.52"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.82%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "295.69"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +7.42%  "
$ws.Range("E7").Value = "  +0.43%  "
$ws.Range("E8").Value = "  -0.24%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.608"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.28%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "43.97"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -5.39%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0927"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.30"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.70%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.90"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.98%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.07"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +22.73%  "
$ws.Range("E15").Value = "  -0.55%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.23"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.16%  "
$ws.Range("D17").Value = "2.587.85"
$ws.Range("E17").Value = "  +0.49%  "
$ws.Range("D18").Value = "2.331.01"
$ws.Range("E18").Value = "  +4.24%  "
$ws.Range("D19").Value = "42.719.65"
$ws.Range("E19").Value = "  -0.25%  "
$ws.Range("E20").Value = "  -0.16%  "
$ws.Range("E21").Value = "  +6.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "75.04"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +4.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.48"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +17.01%  "
$ws.Range("E24").Value = "  +5.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "253.63"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +9.60%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.93"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -3.60%  "
$ws.Range("E27").Value = "  -0.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.55"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -3.16%  "
$ws.Range("E29").Value = "  -0.57%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "37.96"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -5.36%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "175.50"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "22.28"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +5.73%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.16"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -3.49%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0889"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.27%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.68"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +2.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.06"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +9.04%  "
$ws.Range("E37").Value = "  +0.14%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.23"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.94%  "
$ws.Range("E39").Value = "  +1.79%  "
$ws.Range("E40").Value = "  -1.76%  "
$ws.Range("E41").Value = "  -5.62%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "72.00"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.69%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.232"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.39%  "
$ws.Range("E44").Value = "  -0.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.48"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -5.63%  "
$ws.Range("E46").Value = "  -0.52%  "
$ws.Range("E47").Value = "  -2.19%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "106.69"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +6.55%  "
$ws.Range("E49").Value = "  +2.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.64"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +2.34%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.78"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +4.28%  "
